# Add a new row of "mtcars" summary data (row r=14, i.e. the 14 cars with
# 6 cylinders) to the results table. This shifts the existing footer row
# ("Data from the infamous mtcars data set.") down by one row, and updates
# the values of the 4-cylinder / 6-cylinder / 8-cylinder summary rows to
# reflect the new three-group split (cyl=4, cyl=6, cyl=8) used once the
# new row is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the footer row (row 9), pushing the footer
# row (and its merged cell) down to row 10.
$ws.Rows.Item(9).Insert()

# Copy the formatting (styles/borders) of the row above (row 8, a normal
# data row) into the freshly inserted row 9 so the new row's cell styles
# match the existing data rows exactly.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)

# --- Update row 5 (now group with N=1) ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 91
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 2.14
$ws.Range("E5").Value = ""

# --- Update row 6 (now group with N=10) ---
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 81.8
$ws.Range("C6").Value = 21.87235698318771
$ws.Range("D6").Value = 2.3003
$ws.Range("E6").Value = 0.5982073312080948

# --- Update row 7 (group with N=3) ---
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 131.6666666666667
$ws.Range("C7").Value = 37.52776749732568
$ws.Range("D7").Value = 2.755
$ws.Range("E7").Value = 0.1281600561797629

# --- Update row 8 (now group with N=4) ---
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 115.25
$ws.Range("C8").Value = 9.178779875342908
$ws.Range("D8").Value = 3.38875
$ws.Range("E8").Value = 0.1162163929916946

# --- Populate new row 9 (group with N=14) ---
$ws.Range("A9").Value = 14
$ws.Range("B9").Value = 209.2142857142857
$ws.Range("C9").Value = 50.97688551827051
$ws.Range("D9").Value = 3.999214285714287
$ws.Range("E9").Value = 0.7594047444769265
